# Scheduled runner update: refresh cached market-board figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the per-job Leve sheets with freshly pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3145953.5
$ws.Range("J17").Value = 3145953.5
$ws.Range("L17").Value = 9437860.5
$ws.Range("N17").Value = -9438196.5
$ws.Range("H19").Value = 1080.8334
$ws.Range("I19").Value = 495.5
$ws.Range("K19").Value = 495.5
$ws.Range("M19").Value = -320.5
$ws.Range("H34").Value = 3495.2666
$ws.Range("I34").Value = 1409.9231
$ws.Range("J34").Value = 17050
$ws.Range("K34").Value = 1409.9231
$ws.Range("L34").Value = 17050
$ws.Range("M34").Value = -1206.9231
$ws.Range("N34").Value = -17456
$ws.Range("H36").Value = 3495.2666
$ws.Range("I36").Value = 1409.9231
$ws.Range("J36").Value = 17050
$ws.Range("K36").Value = 1409.9231
$ws.Range("L36").Value = 17050
$ws.Range("M36").Value = -694.9231
$ws.Range("N36").Value = -18480
$ws.Range("H62").Value = 15874373
$ws.Range("I62").Value = 23810788
$ws.Range("J62").Value = 1541.4286
$ws.Range("K62").Value = 23810788
$ws.Range("L62").Value = 1541.4286
$ws.Range("M62").Value = -23810164
$ws.Range("N62").Value = -2789.4286
$ws.Range("H65").Value = 15874373
$ws.Range("I65").Value = 23810788
$ws.Range("J65").Value = 1541.4286
$ws.Range("K65").Value = 119053940
$ws.Range("L65").Value = 7707.143
$ws.Range("M65").Value = -119050820
$ws.Range("N65").Value = -13947.143
$ws.Range("H82").Value = 1772.75
$ws.Range("I82").Value = 1772.75
$ws.Range("K82").Value = 5318.25
$ws.Range("M82").Value = -4912.25
$ws.Range("H85").Value = 1772.75
$ws.Range("I85").Value = 1772.75
$ws.Range("K85").Value = 5318.25
$ws.Range("M85").Value = -3914.25
$ws.Range("H116").Value = 8312.058999999999
$ws.Range("I116").Value = 14175.625
$ws.Range("J116").Value = 3100
$ws.Range("K116").Value = 14175.625
$ws.Range("L116").Value = 3100
$ws.Range("M116").Value = -10733.625
$ws.Range("N116").Value = -9984
$ws.Range("H132").Value = 6803599
$ws.Range("I132").Value = 633.5238000000001
$ws.Range("K132").Value = 1900.5714
$ws.Range("M132").Value = 629.4285999999997
$ws.Range("H138").Value = 4613.8525
$ws.Range("I138").Value = 1052
$ws.Range("J138").Value = 7636.0303
$ws.Range("K138").Value = 3156
$ws.Range("L138").Value = 22908.0909
$ws.Range("M138").Value = 1984
$ws.Range("N138").Value = -33188.0909

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4454.642
$ws.Range("I32").Value = 3229.077
$ws.Range("J32").Value = 10077.823
$ws.Range("K32").Value = 3229.077
$ws.Range("L32").Value = 10077.823
$ws.Range("M32").Value = -2942.077
$ws.Range("N32").Value = -10651.823
$ws.Range("H110").Value = 1620.9
$ws.Range("I110").Value = 944.1429000000001
$ws.Range("K110").Value = 944.1429000000001
$ws.Range("M110").Value = 1100.8571
$ws.Range("H122").Value = 1976365.2
$ws.Range("I122").Value = 2568784.8
$ws.Range("K122").Value = 7706354.399999999
$ws.Range("M122").Value = -7703904.399999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1791.804
$ws.Range("I134").Value = 1558.4166
$ws.Range("J134").Value = 2351.9333
$ws.Range("K134").Value = 4675.2498
$ws.Range("L134").Value = 7055.7999
$ws.Range("M134").Value = -2140.2498
$ws.Range("N134").Value = -12125.7999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 9616984
$ws.Range("I16").Value = 19231842
$ws.Range("J16").Value = 2125
$ws.Range("K16").Value = 19231842
$ws.Range("L16").Value = 2125
$ws.Range("M16").Value = -19231555
$ws.Range("N16").Value = -2699
$ws.Range("H31").Value = 8069159.5
$ws.Range("I31").Value = 2105.5518
$ws.Range("J31").Value = 15158389
$ws.Range("K31").Value = 2105.5518
$ws.Range("L31").Value = 15158389
$ws.Range("M31").Value = -1810.5518
$ws.Range("N31").Value = -15158979
$ws.Range("H34").Value = 8069159.5
$ws.Range("I34").Value = 2105.5518
$ws.Range("J34").Value = 15158389
$ws.Range("K34").Value = 2105.5518
$ws.Range("L34").Value = 15158389
$ws.Range("M34").Value = -1903.5518
$ws.Range("N34").Value = -15158793
$ws.Range("H99").Value = 6652.7407
$ws.Range("I99").Value = 6566.1177
$ws.Range("J99").Value = 6800
$ws.Range("K99").Value = 6566.1177
$ws.Range("L99").Value = 6800
$ws.Range("M99").Value = -5068.1177
$ws.Range("N99").Value = -9796
$ws.Range("H113").Value = 9616984
$ws.Range("I113").Value = 19231842
$ws.Range("J113").Value = 2125
$ws.Range("K113").Value = 19231842
$ws.Range("L113").Value = 2125
$ws.Range("M113").Value = -19229672
$ws.Range("N113").Value = -6465
$ws.Range("H126").Value = 6652.7407
$ws.Range("I126").Value = 6566.1177
$ws.Range("J126").Value = 6800
$ws.Range("K126").Value = 19698.3531
$ws.Range("L126").Value = 20400
$ws.Range("M126").Value = -17228.3531
$ws.Range("N126").Value = -25340
$ws.Range("H132").Value = 5716847.5
$ws.Range("I132").Value = 9093106
$ws.Range("J132").Value = 3179.077
$ws.Range("K132").Value = 27279318
$ws.Range("L132").Value = 9537.231
$ws.Range("M132").Value = -27276788
$ws.Range("N132").Value = -14597.231
$ws.Range("H134").Value = 5651605.5
$ws.Range("I134").Value = 8132185.5
$ws.Range("J134").Value = 1395.8889
$ws.Range("K134").Value = 24396556.5
$ws.Range("L134").Value = 4187.6667
$ws.Range("M134").Value = -24394021.5
$ws.Range("N134").Value = -9257.6667
$ws.Range("H141").Value = 29957.684
$ws.Range("J141").Value = 29957.684
$ws.Range("L141").Value = 29957.684
$ws.Range("N141").Value = -40317.684

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 250658.25
$ws.Range("I5").Value = 222.92308
$ws.Range("J5").Value = 546627.25
$ws.Range("K5").Value = 668.76924
$ws.Range("L5").Value = 1639881.75
$ws.Range("M5").Value = -556.76924
$ws.Range("N5").Value = -1640105.75
$ws.Range("H69").Value = 1812.3636
$ws.Range("I69").Value = 500
$ws.Range("J69").Value = 1943.6
$ws.Range("K69").Value = 1500
$ws.Range("L69").Value = 5830.799999999999
$ws.Range("M69").Value = -689
$ws.Range("N69").Value = -7452.799999999999
$ws.Range("H72").Value = 1812.3636
$ws.Range("I72").Value = 500
$ws.Range("J72").Value = 1943.6
$ws.Range("K72").Value = 4500
$ws.Range("L72").Value = 17492.4
$ws.Range("M72").Value = -444
$ws.Range("N72").Value = -25604.4
$ws.Range("H117").Value = 20844920
$ws.Range("J117").Value = 33341692
$ws.Range("L117").Value = 100025076
$ws.Range("N117").Value = -100031960
$ws.Range("H135").Value = 250658.25
$ws.Range("I135").Value = 222.92308
$ws.Range("J135").Value = 546627.25
$ws.Range("K135").Value = 2006.30772
$ws.Range("L135").Value = 4919645.25
$ws.Range("M135").Value = 528.69228
$ws.Range("N135").Value = -4924715.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 47620416
$ws.Range("I113").Value = 90910216
$ws.Range("J113").Value = 1629.2
$ws.Range("K113").Value = 90910216
$ws.Range("L113").Value = 1629.2
$ws.Range("M113").Value = -90908046
$ws.Range("N113").Value = -5969.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 11111837
$ws.Range("I46").Value = 27778386
$ws.Range("J46").Value = 804.44446
$ws.Range("K46").Value = 27778386
$ws.Range("L46").Value = 804.44446
$ws.Range("M46").Value = -27778198
$ws.Range("N46").Value = -1180.44446
$ws.Range("H93").Value = 41667584
$ws.Range("H98").Value = 37500
$ws.Range("J98").Value = 37500
$ws.Range("L98").Value = 37500
$ws.Range("N98").Value = -43490
$ws.Range("H108").Value = 129292.664
$ws.Range("J108").Value = 129292.664
$ws.Range("L108").Value = 129292.664
$ws.Range("N108").Value = -136972.664
$ws.Range("H136").Value = 4506.354
$ws.Range("I136").Value = 2128.6562
$ws.Range("J136").Value = 9261.75
$ws.Range("K136").Value = 6385.9686
$ws.Range("L136").Value = 27785.25
$ws.Range("M136").Value = -3835.9686
$ws.Range("N136").Value = -32885.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2800
$ws.Range("I62").Value = 2800
$ws.Range("K62").Value = 2800
$ws.Range("M62").Value = -2176
$ws.Range("H65").Value = 2800
$ws.Range("I65").Value = 2800
$ws.Range("K65").Value = 14000
$ws.Range("M65").Value = -10880
$ws.Range("H139").Value = 27680
$ws.Range("J139").Value = 27680
$ws.Range("L139").Value = 27680
$ws.Range("N139").Value = -37960
